$p = $ppt.ActivePresentation
# Touching the NotesMaster forces PowerPoint to materialize the Notes
# Master part (and its own theme, saved as ppt/theme/theme2.xml) even
# though the deck has no speaker notes yet.
$nm = $p.NotesMaster
